$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 5
$ws.Range("B3").Value = 4
$ws.Range("B4").Value = 3
$ws.Range("B5").Value = 2
$ws.Range("B6").Value = 2
$ws.Range("B7").Value = 1

$ws.Range("B8").Select()
